$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A207").Value = "IMX-USD"
$ws.Range("A208").Value = "TAO-USD"
$ws.Range("A209").Value = "GRT-USD"
